# Update the marksheet's correct/total marks summary row.
# "Marking" row (B11): points awarded per correct answer -> 5
# "Total" row (B12): total score achieved -> 105
# "Total" row (E12): "correct/total" label text -> "105/140"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 105
$ws.Range("E12").Value = "105/140"
